$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15 ("data_subsets") to hold the new
# "grouping_orders" setting. This shifts rows 15-30 down to 16-31, and
# automatically adjusts the merged cell ranges below it.
$ws.Rows.Item(15).Insert()

# Fill in the new setting name and its comment (value column intentionally
# left blank, matching the template for a not-yet-configured setting).
$ws.Range("A15").Value = "grouping_orders"
$ws.Range("C15").Value = "Control order of groups. Input group names separated by comma. For multiple grouping columns - separate by semicolon."

# Copy the cell formatting (borders/alignment/wrap) from the row below, which
# already carries the standard "setting name" / "comment" styles used
# throughout the table, so the new row matches the table's look.
$ws.Range("A16").Copy()
$ws.Range("A15").PasteSpecial(-4122)

$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# No value for this setting yet, remove the blank cell Insert() created there.
$ws.Range("B15").Clear()

# Match the row height used for similar wrapped, multi-line comment rows.
$ws.Rows.Item(15).RowHeight = 45

# Update the view so the newly added row is visible and selected, same as
# when the author was editing it.
$ws.Range("C15").Select()
$excel.ActiveWindow.ScrollRow = 14
